$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 414.375
$ws.Range("I6").Value = 414.375
$ws.Range("K6").Value = 1243.125
$ws.Range("M6").Value = -1131.125
$ws.Range("H19").Value = 2473.8
$ws.Range("J19").Value = 2592.25
$ws.Range("L19").Value = 2592.25
$ws.Range("N19").Value = -2942.25
$ws.Range("H76").Value = 3805
$ws.Range("I76").Value = 3783.3333
$ws.Range("K76").Value = 3783.3333
$ws.Range("M76").Value = -3468.3333
$ws.Range("H79").Value = 3805
$ws.Range("I79").Value = 3783.3333
$ws.Range("K79").Value = 3783.3333
$ws.Range("M79").Value = -2691.3333
$ws.Range("H92").Value = 445.35715
$ws.Range("I92").Value = 420.07693
$ws.Range("K92").Value = 420.07693
$ws.Range("M92").Value = 827.9230700000001
$ws.Range("H95").Value = 42578.4
$ws.Range("I95").Value = 40567
$ws.Range("J95").Value = 50624
$ws.Range("K95").Value = 40567
$ws.Range("L95").Value = 50624
$ws.Range("M95").Value = -37821
$ws.Range("N95").Value = -56116
$ws.Range("H96").Value = 2749.6667
$ws.Range("I96").Value = 583.8570999999999
$ws.Range("J96").Value = 10330
$ws.Range("K96").Value = 1751.5713
$ws.Range("L96").Value = 30990
$ws.Range("M96").Value = -378.5712999999998
$ws.Range("N96").Value = -33736
$ws.Range("H100").Value = 4819.6
$ws.Range("I100").Value = 5000
$ws.Range("K100").Value = 5000
$ws.Range("M100").Value = -4459
$ws.Range("H106").Value = 2113.7
$ws.Range("I106").Value = 3328.2
$ws.Range("J106").Value = 899.2
$ws.Range("K106").Value = 3328.2
$ws.Range("L106").Value = 899.2
$ws.Range("M106").Value = -2697.2
$ws.Range("N106").Value = -2161.2
$ws.Range("H112").Value = 942.85187
$ws.Range("J112").Value = 995.1739
$ws.Range("L112").Value = 2985.5217
$ws.Range("N112").Value = -5201.5217
$ws.Range("H127").Value = 950
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 950
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 2850
$ws.Range("M127").ClearContents()
$ws.Range("N127").Value = -12770
$ws.Range("H132").Value = 1527.2433
$ws.Range("I132").Value = 1291.9722
$ws.Range("K132").Value = 3875.9166
$ws.Range("M132").Value = -1345.9166
$ws.Range("H137").Value = 324117.5
$ws.Range("I137").Value = 1837.0869
$ws.Range("J137").Value = 661047.0600000001
$ws.Range("K137").Value = 5511.2607
$ws.Range("L137").Value = 1983141.18
$ws.Range("M137").Value = -2961.2607
$ws.Range("N137").Value = -1988241.18
$ws.Range("H138").Value = 2562
$ws.Range("J138").Value = 2443.3333
$ws.Range("L138").Value = 7329.999899999999
$ws.Range("N138").Value = -17609.9999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 13892090
$ws.Range("I45").Value = 3468.6667
$ws.Range("K45").Value = 3468.6667
$ws.Range("M45").Value = -3091.6667
$ws.Range("H102").Value = 101132.37
$ws.Range("I102").Value = 143939
$ws.Range("J102").Value = 26220.75
$ws.Range("K102").Value = 143939
$ws.Range("L102").Value = 26220.75
$ws.Range("M102").Value = -142317
$ws.Range("N102").Value = -29464.75
$ws.Range("H105").Value = 109999.5
$ws.Range("J105").Value = 109999.5
$ws.Range("L105").Value = 109999.5
$ws.Range("N105").Value = -116987.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3894.75
$ws.Range("I94").Value = 3938.611
$ws.Range("K94").Value = 3938.611
$ws.Range("M94").Value = -3487.611
$ws.Range("H132").Value = 30701.701
$ws.Range("J132").Value = 30701.701
$ws.Range("L132").Value = 30701.701
$ws.Range("N132").Value = -40821.701
$ws.Range("H134").Value = 2917.2979
$ws.Range("I134").Value = 2339.8372
$ws.Range("K134").Value = 7019.5116
$ws.Range("M134").Value = -4484.5116

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 6999.143
$ws.Range("I23").Value = 1000
$ws.Range("J23").Value = 7999
$ws.Range("K23").Value = 1000
$ws.Range("L23").Value = 7999
$ws.Range("M23").Value = -760
$ws.Range("N23").Value = -8479
$ws.Range("H27").Value = 6999.143
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 7999
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 7999
$ws.Range("M27").Value = -808
$ws.Range("N27").Value = -8383
$ws.Range("H31").Value = 2541.6667
$ws.Range("J31").Value = 4543.6665
$ws.Range("L31").Value = 4543.6665
$ws.Range("N31").Value = -5133.6665
$ws.Range("H34").Value = 2541.6667
$ws.Range("J34").Value = 4543.6665
$ws.Range("L34").Value = 4543.6665
$ws.Range("N34").Value = -4947.6665
$ws.Range("H62").Value = 2862.125
$ws.Range("I62").Value = 2699.5715
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 2699.5715
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = -2075.5715
$ws.Range("N62").Value = -5248
$ws.Range("H65").Value = 2862.125
$ws.Range("I65").Value = 2699.5715
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 13497.8575
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = -10377.8575
$ws.Range("N65").Value = -26240
$ws.Range("H86").Value = 6681.5454
$ws.Range("I86").Value = 5240
$ws.Range("K86").Value = 5240
$ws.Range("M86").Value = -4117
$ws.Range("H89").Value = 6681.5454
$ws.Range("I89").Value = 5240
$ws.Range("K89").Value = 26200
$ws.Range("M89").Value = -20584
$ws.Range("H134").Value = 2732.111
$ws.Range("I134").Value = 2617.923
$ws.Range("J134").Value = 3029
$ws.Range("K134").Value = 7853.768999999999
$ws.Range("L134").Value = 9087
$ws.Range("M134").Value = -5318.768999999999
$ws.Range("N134").Value = -14157

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 2780
$ws.Range("J54").Value = 3750
$ws.Range("L54").Value = 11250
$ws.Range("N54").Value = -12368
$ws.Range("H57").Value = 3995
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H132").Value = 5088.9
$ws.Range("I132").Value = 1278
$ws.Range("J132").Value = 8899.799999999999
$ws.Range("K132").Value = 11502
$ws.Range("L132").Value = 80098.2
$ws.Range("M132").Value = -8972
$ws.Range("N132").Value = -85158.2
$ws.Range("H139").Value = 2509
$ws.Range("I139").Value = 1338.8125
$ws.Range("K139").Value = 4016.4375
$ws.Range("M139").Value = 1123.5625

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 7000
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H27").Value = 1800
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 1800
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 1800
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -2132
$ws.Range("H107").Value = 635.75757
$ws.Range("I107").Value = 573.375
$ws.Range("J107").Value = 694.4706
$ws.Range("K107").Value = 573.375
$ws.Range("L107").Value = 694.4706
$ws.Range("M107").Value = 1346.625
$ws.Range("N107").Value = -4534.4706
$ws.Range("H113").Value = 3829393.5
$ws.Range("I113").Value = 279802.75
$ws.Range("J113").Value = 6669066
$ws.Range("K113").Value = 279802.75
$ws.Range("L113").Value = 6669066
$ws.Range("M113").Value = -277632.75
$ws.Range("N113").Value = -6673406
$ws.Range("H122").Value = 512755.6
$ws.Range("I122").Value = 722222.5600000001
$ws.Range("K122").Value = 2166667.68
$ws.Range("M122").Value = -2164217.68
$ws.Range("H126").Value = 4503
$ws.Range("I126").Value = 2439
$ws.Range("J126").Value = 6004.091
$ws.Range("K126").Value = 7317
$ws.Range("L126").Value = 18012.273
$ws.Range("M126").Value = -4847
$ws.Range("N126").Value = -22952.273

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4844.154
$ws.Range("I16").Value = 4361.273
$ws.Range("K16").Value = 4361.273
$ws.Range("M16").Value = -4191.273
$ws.Range("H22").Value = 1210.625
$ws.Range("I22").Value = 1287.6
$ws.Range("J22").Value = 1082.3334
$ws.Range("K22").Value = 1287.6
$ws.Range("L22").Value = 1082.3334
$ws.Range("M22").Value = -992.5999999999999
$ws.Range("N22").Value = -1672.3334
$ws.Range("H27").Value = 1210.625
$ws.Range("I27").Value = 1287.6
$ws.Range("J27").Value = 1082.3334
$ws.Range("K27").Value = 1287.6
$ws.Range("L27").Value = 1082.3334
$ws.Range("M27").Value = -1180.6
$ws.Range("N27").Value = -1296.3334
$ws.Range("H46").Value = 2296.4211
$ws.Range("I46").Value = 1544.3334
$ws.Range("K46").Value = 1544.3334
$ws.Range("M46").Value = -1356.3334
$ws.Range("H82").Value = 2682.8333
$ws.Range("J82").Value = 2300
$ws.Range("L82").Value = 2300
$ws.Range("N82").Value = -3022
$ws.Range("H85").Value = 2682.8333
$ws.Range("J85").Value = 2300
$ws.Range("L85").Value = 2300
$ws.Range("N85").Value = -4796
$ws.Range("H93").Value = 2956.5833
$ws.Range("I93").Value = 1813.6666
$ws.Range("K93").Value = 1813.6666
$ws.Range("M93").Value = -565.6666

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 7522717.5
$ws.Range("I96").Value = 1489.8
$ws.Range("K96").Value = 1489.8
$ws.Range("M96").Value = -116.8
